$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# LBCB 1 & 2 shunt calibrations
#
# Update the three shunt-cal readings (B40:B42 = shunt+, zero, shunt-) on
# each of the six load-cell tabs with the new calibration voltages. Every
# other cell on those tabs (C/D/E40:42, M45:M49, and - on the Z2/Z3 tabs -
# A51/A53) is formula-driven off B40:B42 and recalculates automatically.
# ---------------------------------------------------------------------------

$ws_x1 = $wb.Worksheets.Item("X1 LOAD CELL")
$ws_x1.Range("B40").Value = 3.02902
$ws_x1.Range("B41").Value = 0.12922600000000001
$ws_x1.Range("B42").Value = -2.76335

$ws_x2 = $wb.Worksheets.Item("X2 LOAD CELL")
$ws_x2.Range("B40").Value = 2.9851299999999998
$ws_x2.Range("B41").Value = 0.11280800000000001
$ws_x2.Range("B42").Value = -2.7593399999999999

$ws_y1 = $wb.Worksheets.Item("Y1 LOAD CELL")
$ws_y1.Range("B40").Value = 2.9815900000000002
$ws_y1.Range("B41").Value = 0.098238000000000006
$ws_y1.Range("B42").Value = -2.76844

$ws_z1 = $wb.Worksheets.Item("Z1 LOAD CELL")
$ws_z1.Range("B40").Value = 3.4376099999999998
$ws_z1.Range("B41").Value = 0.50129900000000005
$ws_z1.Range("B42").Value = -2.4350800000000001

$ws_z2 = $wb.Worksheets.Item("Z2 LOAD CELL")
$ws_z2.Range("B40").Value = 3.1833999999999998
$ws_z2.Range("B41").Value = 0.277055
$ws_z2.Range("B42").Value = -2.62683

$ws_z3 = $wb.Worksheets.Item("Z3 LOAD CELL")
$ws_z3.Range("B40").Value = 3.2374200000000002
$ws_z3.Range("B41").Value = 0.32539499999999999
$ws_z3.Range("B42").Value = -2.5872000000000002

# Each tab's selection now rests on the updated "zero load" shunt cell.
[void]$ws_x1.Range("B41").Select()
[void]$ws_x2.Range("B41").Select()
[void]$ws_y1.Range("B41").Select()
[void]$ws_z1.Range("B41").Select()
[void]$ws_z2.Range("B41").Select()
[void]$ws_z3.Range("B41").Select()

# ---------------------------------------------------------------------------
# Serial Numbers tab: add a new column C with the freshly-calibrated slope
# (330/"source-point-voltage") for each load cell, and widen column C to fit.
# ---------------------------------------------------------------------------

$ws_sn = $wb.Worksheets.Item("Serial Numbers")
$ws_sn.Range("C3").Value = 32.806390958932603
$ws_sn.Range("C4").Value = 33.309124132459601
$ws_sn.Range("C5").Value = 33.366599678886999
$ws_sn.Range("C6").Value = 32.637019618164203
$ws_sn.Range("C7").Value = 32.982133385550199
$ws_sn.Range("C8").Value = 32.800186152233898

$ws_sn.Columns.Item(3).ColumnWidth = 30.8776

[void]$ws_sn.Range("C8").Select()

# Restore the originally active tab (Z3 LOAD CELL).
[void]$ws_z3.Activate()
